$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (cardholder name / card number) ---
$ws.Range("C2").Value = "Hartmut"
# Force text so the long card number isn't coerced into a (lossy,
# scientific-notation-displayed) number by Excel's General format.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 13.10.2024"

# --- Row 6 (existing transaction, values updated) ---
$ws.Range("B6").Value = "17.10."
$ws.Range("C6").Value = "18.10."
$ws.Range("D6").Value = "ZALANDO MKTPLC EU WEPUFB"
$ws.Range("E6").Value = "116,35-"

# --- Row 7 (existing transaction, values updated) ---
$ws.Range("B7").Value = "21.10."
$ws.Range("C7").Value = "22.10."
$ws.Range("D7").Value = "KARTENZ./21.10 LIDL RO"
$ws.Range("E7").Value = "26,80-"

# --- Row 8 (existing transaction, values updated) ---
$ws.Range("B8").Value = "23.10."
$ws.Range("C8").Value = "24.10."
$ws.Range("D8").Value = "ZALANDO MKTPLC EU HNTYNG"
$ws.Range("E8").Value = "135,28-"

# --- Row 9 (was blank -> new transaction) ---
$ws.Range("B9").Value = "25.10."
$ws.Range("C9").Value = "26.10."
$ws.Range("D9").Value = "KARTENZ./25.10 ALDI SUED RO"
$ws.Range("E9").Value = "121,23-"

# --- Row 10 (was blank -> new transaction) ---
$ws.Range("B10").Value = "29.10."
$ws.Range("C10").Value = "30.10."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-99783751"
$ws.Range("E10").Value = "56,40-"

# --- Row 11 (was blank -> new transaction) ---
$ws.Range("B11").Value = "02.11."
$ws.Range("C11").Value = "03.11."
$ws.Range("D11").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E11").Value = "51,02-"

# E9:E11 previously carried the centered/"total row" styles (s=13 / s=12).
# The diff switches them to the same right-aligned amount style already
# used by E6:E8 (s=17) -- copy that formatting across instead of touching
# each alignment property individually (avoids leaving stray unused style
# records behind).
$ws.Range("E6").Copy()
$ws.Range("E9:E11").PasteSpecial(-4122)
$excel.CutCopyMode = $False

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 07.11.2024"
$ws.Range("E12").Value = "507,08-"

# --- Next billing date line ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 16.11.2024"
